$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: bump Taxonsorteringsordning (B) ---
$ws.Range("B10").Value = 91829

# --- Row 11: bump Taxonsorteringsordning (B) ---
$ws.Range("B11").Value = 80349

# --- Rows 12 & 13: swap the "Järpe" / "Garnlav" sighting records ---

# First drop the cells that must not survive in their new row position.
$ws.Range("AC12").ClearContents()
$ws.Range("L12").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("AF13").ClearContents()
$ws.Range("J13").ClearContents()

# Row 12 becomes the "Garnlav" record (previously row 13's content).
$ws.Range("A12").Value = 131196451
$ws.Range("B12").Value = 79244
$ws.Range("E12").Value = 6425
$ws.Range("F12").Value = "Garnlav"
$ws.Range("G12").Value = "Alectoria sarmentosa"
$ws.Range("H12").Value = "(Ach.) Ach."
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = ""
$ws.Range("N12").Value = ""
$ws.Range("Q12").Value = 500318
$ws.Range("R12").Value = 7016201

# Row 13 becomes the "Järpe" record (previously row 12's content).
$ws.Range("A13").Value = 131196449
$ws.Range("B13").Value = 57064
$ws.Range("E13").Value = 102612
$ws.Range("F13").Value = "Järpe"
$ws.Range("G13").Value = "Tetrastes bonasia"
$ws.Range("H13").Value = "(Linnaeus, 1758)"
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = "2"
$ws.Range("L13").Value = ""
$ws.Range("M13").Value = "födosökande"
$ws.Range("N13").Value = "observerad"
$ws.Range("Q13").Value = 500203
$ws.Range("R13").Value = 7016330
$ws.Range("AC13").Value = "Synobservation av 2 st födosökande järpar."

# --- Rows 14 & 15: swap the "Talltita" / "Garnlav" sighting records ---

$ws.Range("AC14").ClearContents()
$ws.Range("L14").ClearContents()
$ws.Range("M14").ClearContents()
$ws.Range("AF15").ClearContents()
$ws.Range("J15").ClearContents()

# Row 14 becomes the "Garnlav" record (previously row 15's content).
$ws.Range("A14").Value = 131196452
$ws.Range("B14").Value = 79244
$ws.Range("E14").Value = 6425
$ws.Range("F14").Value = "Garnlav"
$ws.Range("G14").Value = "Alectoria sarmentosa"
$ws.Range("H14").Value = "(Ach.) Ach."
$ws.Range("I14").Value = ""
$ws.Range("J14").Value = ""
$ws.Range("N14").Value = ""
$ws.Range("Q14").Value = 500345
$ws.Range("R14").Value = 7016371

# Row 15 becomes the "Talltita" record (previously row 14's content).
$ws.Range("A15").Value = 131196447
$ws.Range("B15").Value = 58043
$ws.Range("E15").Value = 103021
$ws.Range("F15").Value = "Talltita"
$ws.Range("G15").Value = "Poecile montanus"
$ws.Range("H15").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = "1"
$ws.Range("L15").Value = ""
$ws.Range("M15").Value = "förbiflygande"
$ws.Range("N15").Value = "observerad"
$ws.Range("Q15").Value = 500269
$ws.Range("R15").Value = 7016195
$ws.Range("AC15").Value = "Synobservation av 1 st talltita."
